$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.054.69'
$ws.Range('E2').Value = '  +11.39%  '
$ws.Range('D3').Value = '1.811.65'
$ws.Range('E3').Value = '  +7.86%  '
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').Value = '227.84'
$ws.Range('E5').Value = '  +3.56%  '
$ws.Range('D6').Value = '0.543'
$ws.Range('E6').Value = '  +3.23%  '
$ws.Range('E7').Value = '  -0.24%  '
$ws.Range('E8').Value = '  +4.42%  '
$ws.Range('D9').Value = '46.57'
$ws.Range('E9').Value = '  +5.31%  '
$ws.Range('E11').Value = '  +6.07%  '
$ws.Range('D12').Value = '0.0928'
$ws.Range('E12').Value = '  +2.16%  '
$ws.Range('D13').Value = '2.070.69'
$ws.Range('E13').Value = '  +7.81%  '
$ws.Range('D14').Value = '1.807.86'
$ws.Range('E14').Value = '  +7.71%  '
$ws.Range('D15').Value = '0.639'
$ws.Range('E15').Value = '  +2.79%  '
$ws.Range('D16').Value = '33.995.86'
$ws.Range('E16').Value = '  +11.18%  '
$ws.Range('D17').Value = '10.18'
$ws.Range('E17').Value = '  -2.61%  '
$ws.Range('E18').Value = '  +7.05%  '
$ws.Range('D19').Value = '69.31'
$ws.Range('E19').Value = '  +4.41%  '
$ws.Range('D20').Value = '256.79'
$ws.Range('E20').Value = '  +4.79%  '
$ws.Range('E21').Value = '  +4.04%  '
$ws.Range('E22').Value = '  -0.17%  '
$ws.Range('E23').Value = '  +3.12%  '
$ws.Range('D24').Value = '4.33'
$ws.Range('E24').Value = '  +1.10%  '
$ws.Range('E25').Value = '  +1.97%  '
$ws.Range('D26').Value = '158.53'
$ws.Range('E26').Value = '  +0.57%  '
$ws.Range('D27').Value = '16.57'
$ws.Range('E27').Value = '  +4.35%  '
$ws.Range('D28').Value = '7.11'
$ws.Range('E28').Value = '  +5.86%  '
$ws.Range('E29').Value = '  +2.35%  '
$ws.Range('E30').Value = '  -0.17%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').Value = '3.86'
$ws.Range('E31').Value = '  +11.01%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').Value = '0.0512'
$ws.Range('E32').Value = '  +2.88%  '
$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('E33').Value = '  +4.73%  '
$ws.Range('B34').Value = 'MinaProtocolToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/J7st_qGwz+minaprotocoltoken-mina'
$ws.Range('D34').Value = '1.65'
$ws.Range('E34').Value = '  +300.38%  '
$ws.Range('D35').Value = '3.49'
$ws.Range('E35').Value = '  +5.98%  '
$ws.Range('D36').Value = '1.536.82'
$ws.Range('E36').Value = '  +1.79%  '
$ws.Range('E37').Value = '  +2.25%  '
$ws.Range('E38').Value = '  +4.06%  '
$ws.Range('D39').Value = '84.27'
$ws.Range('E39').Value = '  +0.54%  '
$ws.Range('E40').Value = '  +4.70%  '
$ws.Range('D41').Value = '0.621'
$ws.Range('E41').Value = '  +5.03%  '
$ws.Range('E42').Value = '  +3.24%  '
$ws.Range('E43').Value = '  +1.57%  '
$ws.Range('E44').Value = '  +8.19%  '
$ws.Range('E45').Value = '  +6.90%  '
$ws.Range('E46').Value = '  +3.80%  '
$ws.Range('E47').Value = '  +4.00%  '
$ws.Range('D48').Value = '1.968.76'
$ws.Range('E48').Value = '  +8.29%  '
$ws.Range('D49').Value = '5.72'
$ws.Range('E49').Value = '  +2.64%  '
$ws.Range('E50').Value = '  -0.28%  '
$ws.Range('D51').Value = '52.44'
$ws.Range('E51').Value = '  +1.88%  '

# Trailing-zero numeric-looking strings: force text to avoid Excel number coercion
$orig_D33 = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.20'
$ws.Range('D33').Style = $orig_D33
$orig_D37 = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.80'
$ws.Range('D37').Style = $orig_D37
$orig_D42 = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.80'
$ws.Range('D42').Style = $orig_D42
